$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.048.43'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.661.58'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -4.06%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.55'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.18'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -6.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.655.74'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.13%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.522'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.14'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.461'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.67'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.85%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -6.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.274.61'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.659.54'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.084.25'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -4.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.10'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -6.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.81'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '489.68'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -4.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.06'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -6.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.715'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.09'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.68%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -7.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000139'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.16'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -5.32%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.94'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -6.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.91'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.36'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -6.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.73'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.73'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.803.61'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.107'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -6.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.600.81'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.988'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.75'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.27%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -7.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.321'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '440.58'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -8.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '48.47'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.65%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -7.48%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -8.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.31'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.36%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '141.42'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '39.61'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -10.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.747.47'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -7.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0346'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.97%  '
